$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are kept as plain text, matching the
# original "inline string" representation, rather than being auto-
# converted to numbers by Excel when they look numeric.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.788.65"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.303.34"
$ws.Range("E3").Value = "  +5.84%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.39"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.31"
$ws.Range("E6").Value = "  +4.42%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.301.52"
$ws.Range("E8").Value = "  +5.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("E10").Value = "  +2.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.50"
$ws.Range("E11").Value = "  +4.63%  "
$ws.Range("E12").Value = "  +3.47%  "
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.70"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.848.67"
$ws.Range("E15").Value = "  +5.91%  "
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.308.08"
$ws.Range("E17").Value = "  +6.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.886.38"
$ws.Range("E18").Value = "  +1.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.88"
$ws.Range("E19").Value = "  +3.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.03"
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("E22").Value = "  +5.13%  "
$ws.Range("E23").Value = "  +4.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.83"
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.44"
$ws.Range("E25").Value = "  +4.14%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.36"
$ws.Range("E28").Value = "  +6.56%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.14"
$ws.Range("E30").Value = "  +3.28%  "
$ws.Range("E31").Value = "  +4.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.28"
$ws.Range("E32").Value = "  +9.60%  "
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.53"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.10"
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.89"
$ws.Range("E37").Value = "  +1.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0749"
$ws.Range("E38").Value = "  +7.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0402"
$ws.Range("E39").Value = "  +4.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "429.02"
$ws.Range("E40").Value = "  +2.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.046.41"
$ws.Range("E41").Value = "  +5.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.40"
$ws.Range("E42").Value = "  +2.20%  "
$ws.Range("E43").Value = "  +2.71%  "
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.266"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("E46").Value = "  +4.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.29"
$ws.Range("E47").Value = "  +3.29%  "
$ws.Range("E48").Value = "  +0.01%  "

# Row 49: Arweave -> Stellar
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.115"
$ws.Range("E49").Value = "  +2.28%  "

# Row 50: Stellar -> ThetaToken
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.31"
$ws.Range("E50").Value = "  +2.70%  "

# Row 51: ThetaToken -> Monero
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.27"
$ws.Range("E51").Value = "  +3.79%  "
